$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (D: Price, E: Volume(1h), and B/C for the Maker/VeChain row swap).
# NumberFormat is forced to text ("@") before assignment so values such as "1.003" or
# "0.6278" are stored as literal text (matching the source inlineStr cells) instead of
# being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.103.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6278"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07486"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2925"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07685"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.827.41"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.006"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6672"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.69"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009345"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -8.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.984"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.109.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.081.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "223.12"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.096"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.85"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1389"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.483"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05725"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.147"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.077"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.207"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7418"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.830"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.671"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01778"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.212.83"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.517"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8898"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.979.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.48"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5098"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4058"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.009"
